$d = $word.ActiveDocument

# The lone paragraph currently reads "16.04.2010 - " + "USB Host hatte zu
# wenig Strom geliefert offenbar (UART geht/geht nicht)" as two runs, with
# the "_GoBack" bookmark sitting between them. The target document instead
# has that sentence as a single run in paragraph 1, a blank paragraph 2, and
# a new paragraph 3 ("ARP, TCP,UDP,IP ...") that ends with the bookmark.

# Drop the old bookmark; it gets rebuilt at the end of the new paragraph 3.
$d.Bookmarks("_GoBack").Delete()

# Empty out paragraph 1 completely first. Doing this turns the existing,
# already-in-the-file paragraph into the blank middle paragraph later on
# (keeping it a clean, run-less paragraph), while the sentence and the new
# heading are written into brand-new paragraphs inserted around it.
$p1 = $d.Paragraphs(1)
$rClear = $p1.Range
$rClear.MoveEnd(1, -1) | Out-Null
$rClear.Text = ""

# Insert a new (still empty) paragraph before the blank one, and another
# new (still empty) paragraph after it - this is where the sentence and the
# heading will go, leaving the original paragraph object as the blank one
# sandwiched in between.
$p1.Range.InsertParagraphBefore()
$pBlank = $d.Paragraphs(2)
$pBlank.Range.InsertParagraphAfter()

# Fill in paragraph 1 with the merged sentence (single run).
$pFirst = $d.Paragraphs(1)
$rFirst = $pFirst.Range
$rFirst.MoveEnd(1, -1) | Out-Null
$rFirst.Text = "16.04.2010 - USB Host hatte zu wenig Strom geliefert offenbar (UART geht/geht nicht)"

# Fill in paragraph 3 with the new heading text, plus a one-character "X"
# placeholder appended right after it.
$pThird = $d.Paragraphs(3)
$rThird = $pThird.Range
$rThird.MoveEnd(1, -1) | Out-Null
$rThird.Text = "ARP, TCP,UDP,IP – Programmablauf und Funktionsweise des µIPX"

# Re-create the "_GoBack" bookmark as a collapsed point right after the new
# heading text (i.e. before paragraph 3's paragraph mark). Anchoring it while
# the "X" placeholder still follows keeps Bookmarks.Add from snapping the
# position to a neighbouring paragraph; the placeholder is deleted right
# after, leaving the bookmark collapsed exactly at the end of the heading.
$endPos = $pThird.Range.End - 2
$bmRange = $d.Range($endPos, $endPos)
$bmRange.Bookmarks.Add("_GoBack") | Out-Null
$placeholder = $d.Range($endPos, $endPos + 1)
$placeholder.Text = ""
